# Applies the row 4/5/6 data rotation described by the commit diff:
#   - row 4 takes on what used to be row 6's record (minus the "Riklig" comment)
#   - row 5 takes on what used to be row 4's record
#   - row 6 takes on what used to be row 5's record (plus the "Riklig" comment)
#
# Columns A, B, E, Q, R are numeric; the rest (D, F, G, H, I, J, Y, AA, AC) are
# stored as text even when the content looks numeric/date-like, so those are
# written with a leading apostrophe to force Excel to keep them as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 -------------------------------------------------------------
$ws.Cells.Item(4,1).Value  = 111782531          # A4 Id
$ws.Cells.Item(4,2).Value  = 93388               # B4 Taxonsorteringsordning
$ws.Cells.Item(4,4).Value  = "'LC"               # D4 Rödlistade
$ws.Cells.Item(4,5).Value  = 2180                # E4 TaxonId
$ws.Cells.Item(4,6).Value  = "'Blåmossa"         # F4 Artnamn
$ws.Cells.Item(4,7).Value  = "'Leucobryum glaucum"   # G4 Vetenskapligt namn
$ws.Cells.Item(4,8).Value  = "'(Hedw.) Ångstr."   # H4 Auktor
$ws.Cells.Item(4,9).Value  = "'10"               # I4 Antal
$ws.Cells.Item(4,10).Value = "'plantor/tuvor"     # J4 Enhet
$ws.Cells.Item(4,17).Value = 572404.6564225279    # Q4 Ost
$ws.Cells.Item(4,18).Value = 6300359.898093811    # R4 Nord
$ws.Cells.Item(4,29).Value = "'Riklig"            # AC4 Publik kommentar (new)

# --- Row 5 -------------------------------------------------------------
$ws.Cells.Item(5,1).Value  = 111782529
$ws.Cells.Item(5,2).Value  = 88869
$ws.Cells.Item(5,4).Value  = "'NT"
$ws.Cells.Item(5,5).Value  = 2008
$ws.Cells.Item(5,6).Value  = "'Fyrflikig jordstjärna"
$ws.Cells.Item(5,7).Value  = "'Geastrum quadrifidum"
$ws.Cells.Item(5,8).Value  = "'Pers.:Pers."
$ws.Cells.Item(5,10).Value = "'fruktkroppar"
$ws.Cells.Item(5,17).Value = 572410.864430059
$ws.Cells.Item(5,18).Value = 6300349.611970967
$ws.Cells.Item(5,25).Value = "'2022-09-13"        # Y5 Startdatum
$ws.Cells.Item(5,27).Value = "'2022-09-13"        # AA5 Slutdatum

# --- Row 6 -------------------------------------------------------------
$ws.Cells.Item(6,1).Value  = 111782537
$ws.Cells.Item(6,9).Value  = "'3"                 # I6 Antal
$ws.Cells.Item(6,17).Value = 572143.9508974494
$ws.Cells.Item(6,18).Value = 6300252.539621729
$ws.Cells.Item(6,25).Value = "'2022-09-12"        # Y6 Startdatum
$ws.Cells.Item(6,27).Value = "'2022-09-12"        # AA6 Slutdatum
$ws.Cells.Item(6,29).ClearContents()              # AC6 Publik kommentar (removed)
